$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# EMU -> point conversion (1 pt = 12700 EMU)
$left   = 397950   / 12700.0
$top    = 1154050  / 12700.0
$width  = 11023200 / 12700.0
$height = 2986200  / 12700.0

$shape = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shape.Name = "Google Shape;120;p17"

# spPr: no fill, no line
$shape.Fill.Visible = 0
$shape.Line.Visible = 0

$tf = $shape.TextFrame
$tf.WordWrap = 1
$tf.AutoSize = 1
$tf.VerticalAnchor = 1
$tf.MarginLeft = 91425 / 12700.0
$tf.MarginRight = 91425 / 12700.0
$tf.MarginTop = 91425 / 12700.0
$tf.MarginBottom = 91425 / 12700.0

$tr = $tf.TextRange
$tr.Text = "Knapsack (items, weight):"

$tr.Font.Name = "Calibri"
$tr.Font.NameFarEast = "Calibri"
$tr.Font.NameComplexScript = "Calibri"
$tr.ParagraphFormat.Alignment = 1
$tr.IndentLevel = 0
$tr.ParagraphFormat.Bullet.Type = 0
$tr.ParagraphFormat.SpaceBefore = 0
$tr.ParagraphFormat.SpaceAfter = 0

$tr.InsertAfter("`r  Knapsackinator (weights_remain, i_index):")
$tr.InsertAfter("`r    if weights_remain is empty or i_index is greater than or equal to the length of the items:")
$tr.InsertAfter("`r      return a list containing zero and an empty list")
$tr.InsertAfter("`r    with_item <- Knapsackinator (weights_remain - items's weight, i_index + 1)")
$tr.InsertAfter("`r    without_item <- Knapsackinator (weights_remain, i_index + 1)")
$tr.InsertAfter("`r    add item's value to with_item array")
$tr.InsertAfter("`r    check if with_item's value is greater than without_item")
$tr.InsertAfter("`r      return with_item")
$tr.InsertAfter("`r    else:")
$tr.InsertAfter("`r      return without_item")
$tr.InsertAfter("`r  return Knapsackinator(weight, 0)")
$tr.InsertAfter("`r")

Write-Output "New shape added: $($shape.Name) id=$($shape.Id) shapes=$($s.Shapes.Count)"
